# Repull data, push all data, mean calculation
# Updates the "dSF" (column F) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 3
$ws.Range("F7").Value = -13
$ws.Range("F8").Value = -6
$ws.Range("F11").Value = -9
$ws.Range("F12").Value = 13
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -3
